$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

# Fill in the result for match 6 (Santa Cruz x Tropinha) on row 7
$ws.Range("E7").Value = "6x3"
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = "Finalizado"

# Leave selection on K7, matching the author's final cursor position
$ws.Range("K7").Select()
